{"js": "const replacements = [\n  [\"2024-01-17 Wednesday\", \"2024-01-18 Thursday\"],\n  [\"648\u00f79=72, 0\", \"511\u00f72=255, 1\"],\n  [\"474\u00f77=67, 5\", \"384\u00f73=128, 0\"],\n  [\"742\u00f79=82, 4\", \"868\u00f75=173, 3\"],\n  [\"928\u00f78=116, 0\", \"387\u00f72=193, 1\"],\n  [\"361\u00f79=40, 1\", \"701\u00f79=77, 8\"],\n  [\"148\u00f72=74, 0\", \"602\u00f73=200, 2\"],\n  [\"380\u00f75=76, 0\", \"186\u00f73=62, 0\"],\n  [\"215\u00f75=43, 0\", \"389\u00f72=194, 1\"],\n  [\"946\u00f78=118, 2\", \"713\u00f74=178, 1\"],\n  [\"366\u00f72=183, 0\", \"715\u00f77=102, 1\"],\n  [\"204\u00f79=22, 6\", \"101\u00f77=14, 3\"],\n  [\"455\u00f79=50, 5\", \"630\u00f75=126, 0\"],\n  [\"648\u00f76=108, 0\", \"562\u00f76=93, 4\"],\n  [\"879\u00f79=97, 6\", \"238\u00f78=29, 6\"],\n  [\"504\u00f73=168, 0\", \"384\u00f75=76, 4\"],\n  [\"640\u00f72=320, 0\", \"990\u00f72=495, 0\"],\n  [\"474\u00f73=158, 0\", \"737\u00f79=81, 8\"],\n  [\"790\u00f73=263, 1\", \"704\u00f79=78, 2\"],\n  [\"980\u00f75=196, 0\", \"884\u00f79=98, 2\"],\n  [\"913\u00f77=130, 3\", \"125\u00f79=13, 8\"],\n  [\"628\u00f79=69, 7\", \"881\u00f79=97, 8\"],\n  [\"365\u00f78=45, 5\", \"437\u00f72=218, 1\"],\n  [\"917\u00f73=305, 2\", \"260\u00f76=43, 2\"],\n  [\"777\u00f75=155, 2\", \"390\u00f74=97, 2\"],\n  [\"248\u00f76=41, 2\", \"400\u00f74=100, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-17 Wednesday\", \"2024-01-18 Thursday\"),\n    @(\"648\u00f79=72, 0\", \"511\u00f72=255, 1\"),\n    @(\"474\u00f77=67, 5\", \"384\u00f73=128, 0\"),\n    @(\"742\u00f79=82, 4\", \"868\u00f75=173, 3\"),\n    @(\"928\u00f78=116, 0\", \"387\u00f72=193, 1\"),\n    @(\"361\u00f79=40, 1\", \"701\u00f79=77, 8\"),\n    @(\"148\u00f72=74, 0\", \"602\u00f73=200, 2\"),\n    @(\"380\u00f75=76, 0\", \"186\u00f73=62, 0\"),\n    @(\"215\u00f75=43, 0\", \"389\u00f72=194, 1\"),\n    @(\"946\u00f78=118, 2\", \"713\u00f74=178, 1\"),\n    @(\"366\u00f72=183, 0\", \"715\u00f77=102, 1\"),\n    @(\"204\u00f79=22, 6\", \"101\u00f77=14, 3\"),\n    @(\"455\u00f79=50, 5\", \"630\u00f75=126, 0\"),\n    @(\"648\u00f76=108, 0\", \"562\u00f76=93, 4\"),\n    @(\"879\u00f79=97, 6\", \"238\u00f78=29, 6\"),\n    @(\"504\u00f73=168, 0\", \"384\u00f75=76, 4\"),\n    @(\"640\u00f72=320, 0\", \"990\u00f72=495, 0\"),\n    @(\"474\u00f73=158, 0\", \"737\u00f79=81, 8\"),\n    @(\"790\u00f73=263, 1\", \"704\u00f79=78, 2\"),\n    @(\"980\u00f75=196, 0\", \"884\u00f79=98, 2\"),\n    @(\"913\u00f77=130, 3\", \"125\u00f79=13, 8\"),\n    @(\"628\u00f79=69, 7\", \"881\u00f79=97, 8\"),\n    @(\"365\u00f78=45, 5\", \"437\u00f72=218, 1\"),\n    @(\"917\u00f73=305, 2\", \"260\u00f76=43, 2\"),\n    @(\"777\u00f75=155, 2\", \"390\u00f74=97, 2\"),\n    @(\"248\u00f76=41, 2\", \"400\u00f74=100, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $result = $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2, [ref]$true, [ref]$false, [ref]$false, [ref]$false)\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
